$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Table support: add a fourth row holding the "total amount" placeholder,
# directly below the existing ${amount} row.
$ws.Range("A4").Value = '${totalAmount}'

# Zoom in on the (now taller) sheet, matching the saved view state.
$excel.ActiveWindow.Zoom = 385

# Select the newly added cell, matching the saved view state.
$ws.Range("A4").Select()
